# "Common" sheet: add a new "VSTAT License File" row right before the
# existing "SD-WAN Portal License File" row (it becomes the new row 70),
# shifting every row below it down by one.
#
# Row/format/merged-cell/data-validation shifting is handled for free by
# Rows.Insert(), but cell *comments* in this engine are anchored to a
# fixed row index and are NOT moved by Insert(), so they have to be
# relocated by hand (bottom-up, so we never overwrite a comment we still
# need to read) before the new comment for row 70 is created.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Common")

$insertRow = 70
$lastRowBeforeInsert = 129

# Shift rows 70:129 -> 71:130 (values, styles, merged cells, data
# validations all move automatically).
$ws.Rows($insertRow).Insert()

# Re-home every comment that used to sit on row N (70..129) onto row
# N+1, walking from the bottom up so a row is always read before it
# gets overwritten by the row above it moving into it.
for ($r = $lastRowBeforeInsert; $r -ge $insertRow; $r--) {
    $srcCell = $ws.Cells.Item($r, 1)
    $dstCell = $ws.Cells.Item($r + 1, 1)
    if ($srcCell.Comment -ne $null) {
        $commentText = $srcCell.Comment.Text()
        $dstCell.AddComment($commentText) | Out-Null
    } elseif ($dstCell.Comment -ne $null) {
        # Destination holds a stale comment from before the insert
        # (e.g. a section-header row that never had one) - clear it.
        $dstCell.Comment.Delete()
    }
}

# Give the newly inserted row the same look as the other "leaf" rows in
# this block (style ids matching its new neighbour at row 71).
$ws.Range("A71:B71").Copy()
$ws.Range("A70:B70").PasteSpecial(-4122)

$ws.Cells.Item($insertRow, 1).Value = "VSTAT License File"
$ws.Cells.Item($insertRow, 1).AddComment("Optional License File for Elasticsearch [default: ]") | Out-Null
